$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$insertionPoint = $lastPara.Range

# --- Paragraph 0 (ilvl=0) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 1
$nr.Text = 'Code in '
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('a')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('pp_')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('u')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('tils.py')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' lines 131-148')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 1 (ilvl=0) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 1
$nr.Text = 'Explain your choices form question 2. Describe what happens when you use a transformer based encoder model and how this relates/informs to embedding a document for storage in a vector database.'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 2 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'Code is commented'
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' in app_utils')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('.py')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('. ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 3 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'I chose to use RecursiveCharacterTextSpliter because it is the standard text splitter in langchain that splits text and allows for overlap.'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' I chose to overwrite ‘docs’ with the chunked data, because it appears that that is what later code is expecting ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('and')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' we already have the ‘orig_docs’ saved. I chose chunk_size = 1000 because it simply seemed like a reasonable value for this task. I chose overlap = 200 because this is 20% of 1000 (seems to be the standard overlap proportion for similar functions). Some applications seem so use chunk sizes of 5000+, but I’m ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('guessing')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' we')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' would use shorter ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('documents for testing purposes here.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 4 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'Using a transformer-based encoder, the input document(s) is broken into smaller chunks (or tokens) and encoded by the model in a way that preserves context for each token. Smaller chunks mean we’re more likely to be under the token max of our embedding model. '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('Transformer-based embeddings are typically vectors of numeric values (which can represent words and meanings in high-dimensional space), which can be stored in vector databases (unlike the original unembedded data).  V')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('ector database')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('s')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' (like FAISS) ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('allow efficient storage & retrieval of ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('these ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('embedding')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('s')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 5 (ilvl=0) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 1
$nr.Text = 'Code in app_utils.py'
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' lines 91-134')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 6 (ilvl=0) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 1
$nr.Text = 'The code currently assumes the answer is contained in exactly 1 chunk and the vector'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('search (embedding model similarity) will grab the right chunk. Discuss why this is')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('problematic and what type of questions might be difficult to answer.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 7 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'The most obvious problem is if the best answer to a prompt/question '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('(or the necessary context for that answer)')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' is spread across multiple chunks. Relying on only one chunk can lead to inaccuracies in this case. ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 8 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'Vector similarity may not sufficiently identify the best answer/chunk. Maybe a question or answer is ambiguous and requires a very complex understanding of the text beyond one chunk. Maybe the dialect the text was written in is different from the dialect of the question, resulting in decreased vector similarity that is unrelated to sentence meaning.'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 9 (ilvl=0) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 1
$nr.Text = 'Code in app_utils.py'
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' in generate_kb_response function; lines 183 – 223')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 10 (ilvl=0) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 1
$nr.Text = 'Discuss your modification in 6 and the limitations, pros, and cons when determining'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('how much context to bring in')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 11 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'I modified the code to grab our relevant chunks '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('plus')
$nr.Font.Bold = $true
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' the chunks ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('around')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' and ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('between')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' those chunks. The buffer I arbitrarily set to 3 (as in, grab everything from 3 chunks before the first relevant chunk to 3 chunks after the last relevant chunk.).')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 12 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'The Pros are increased context, better text generation, more accurate answers and a more robust model understanding of the text. '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 13 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'The cons are '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('the increased computational cost; additional short-term memory usage to store the additional text chunks (though, I’m assuming we’re working in contexts where this is negligible); potentially exceeding the maximum prompt length; possibly diluting the power of the most relevant chunks by surrounding them with a collection of other less relevant chunks.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 14 (ilvl=0) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 1
$nr.Text = 'Discuss the overall project, what a decoder transformer model is (and how different'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('from encoder), what types of questions the chatbot could answer, and ways one might')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('improve RAG to handle a broader range of questions on your pdf')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 15 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'The overall project takes documents for context and processes them to create an embedding representation of them, and uses them to generate contextualized text responses using any one of multiple text generation APIs.  It can potentially circumvent the input limits of using APIs like ChatGPT directly by storing and '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('processing documents locally and selectively providing a few chunks from a large body of work that are best for generating relevant text responses. ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 16 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'A decoder transformer model '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('“decodes” embeddings to generate text (or other output) in response to some prompt. Good for things like AI chatbots or machine translation. In contrast, an encoder processes input (in this context, text) to encode it into a vector representation. Encoder-decoder models (unlike decoder-only models) can potentially provide contextual representations of text that can be fine-tuneed on decoder tasks. ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 17 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'This app allows you to use AI to explore text documents quickly and easily.  This chatbot could answer questions about the content of the text(s), answer questions based on the content of the text(s), interpret the text(s) for clarification or explaining things to user(s). The chatbot could even give advice for improving the documents, or generating new documents based on the content of the documents uploaded.'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 18 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'RAG (Retrieval-Augmented Generation) '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('can be improved by improving retrieval mechanisms, e.g. using multiple retrieval models and combining them somehow (as with MM-RAG); by fine-tuning your model on a particular domain or dataset related to the kinds of prompts you want to ask. RAG could also be improved by implementing a step in the model where the user prompts are rewritten/paraphrased by an interpreter (prior to text generation) to better match the language in the uploaded documents. It can also be improved by implementing a feedback loop where the model learns from user interactions.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 19 (ilvl=0) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 1
$nr.Text = 'EXTRA CREDIT: Discuss how one might update “./src/app_utils.py” (but don’t update) to also incorporate the chat history. For example, if a user asks a follow-up question that is only understood given the previous LLM response. What complications could arise and how could you handle them?'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 20 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'Update the structure to '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('store the chat history locally and input it in each prompt as additional context.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 21 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'OR Update the structure to store the chat history in a history component of the database, or a queue of previous interactions.'
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 22 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = ' '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('One complication would be ambiguity. The model might not understand how to query the chat history for relevant context. If we’re including chat history as part of the “relevant context” portion of the code, this may also contribute to excessively long prompts, dilution of our relevant docs. You may want to dynamically adjust the context padding we implemented in question 6 to increase how many chunks you’re grabbing for context, or eventually you may wind up with a situation where most of the “most relevant docs” are from the chat history and ')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('not')
$nr.Font.Bold = $false
$nr.Font.Italic = $true
$nr.Font.ItalicBi = $true
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' from the uploaded documents.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

# --- Paragraph 23 (ilvl=1) ---
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$newPara.Range.ListFormat.ListLevelNumber = 2
$nr.Text = 'Another complication would be memory & performance. How would the history be stored? Where? For how long? Storing/processing a large chat history could '
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter('impact RAM & model performance.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$nr.Collapse(0)
$nr.InsertAfter(' It could also be a data privacy/security concern.')
$nr.Font.Bold = $false
$nr.Font.Italic = $false
$nr.Font.ItalicBi = $false
$nr.Font.Size = 12
$insertionPoint = $newPara.Range

Write-Output "Done inserting paragraphs"
Write-Output $d.Paragraphs.Count
